$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Matrix": drop the User/Password columns (D:E) and rename the
# "Nombre" column header/values to Login / Logout.
# ---------------------------------------------------------------------
$matrix = $wb.Worksheets.Item("Matrix")

$matrix.Range("B2").Value = "Login"
$matrix.Range("B3").Value = "Logout"

$matrix.Columns("D:E").Delete()

# ---------------------------------------------------------------------
# Sheet "TC001": replace the placeholder Oscar1..Oscar10 sample data
# with the real admin login test case content.
# ---------------------------------------------------------------------
$tc001 = $wb.Worksheets.Item("TC001")

$tc001.Cells.Clear()

$tc001.Range("A1").Value = "ID"
$tc001.Range("B1").Value = "Test Case Name"
$tc001.Range("E1").Value = "Expected result"
$tc001.Range("D1").Value = "Actual result"
$tc001.Range("C1").Value = "Parameter"

$tc001.Range("A2").Value = "TC001"
$tc001.Range("B2").Value = "Log In as Admin"
$tc001.Range("D2").Value = "Type Authoriacion code"
$tc001.Range("E2").Value = "Authotizarion code should be entered"
$tc001.Range("C2").Value = "test"

$tc001.Range("C3").Value = "test.admin@ur.com"
$tc001.Range("D3").Value = "Type user name into the UserName text field"
$tc001.Range("E3").Value = "User name should be introduces"
$tc001.Range("A3").Value = "TC001"

$tc001.Hyperlinks.Add($tc001.Range("C3"), "mailto:test.admin@ur.com")

$tc001.Columns("A").ColumnWidth = 6
$tc001.Columns("B").ColumnWidth = 14.166666666666666
$tc001.Columns("C").ColumnWidth = 17.833333333333332
$tc001.Columns("D").ColumnWidth = 21.333333333333332
$tc001.Columns("E").ColumnWidth = 34.333333333333336
$tc001.Columns("F").ColumnWidth = 7

# ---------------------------------------------------------------------
# Selections / active sheet: Matrix ends up with B3 selected, TC001
# becomes the visible/active sheet with E3 selected.
# ---------------------------------------------------------------------
$matrix.Range("B3").Select()

$tc001.Activate()
$tc001.Range("E3").Select()
